$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers "zmin" (J1) and "zmax" (K1), matching the centered
# style used by the other header cells (A1:I1).
$ws.Range("J1").Value = "zmin"
$ws.Range("K1").Value = "zmax"
$ws.Range("J1:K1").HorizontalAlignment = -4108

# Fill the zmin/zmax data columns for every data row (2-13).
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 10).Value = 0.2
    $ws.Cells.Item($r, 11).Value = 1
}
